# Achievements_Readme.docx - June 2022 samples publish (#11)
#
# The readme text was re-flowed by Word's proofing pass: a handful of
# paragraphs that used to hold a single run now hold several runs split
# right where the grammar/spell checker flagged a phrase, with
# <w:proofErr/> markers bracketing the flagged text. This script
# reproduces those paragraphs (same paragraph identity / pPr, new run
# split + proofErr markers) using Range.InsertXML, which accepts a
# block-level WordprocessingML fragment and splices it in place of the
# paragraph's current content.

$d = $word.ActiveDocument

function Get-ParaByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

function Replace-ParaXml($needle, $xml) {
    $p = Get-ParaByText($needle)
    if ($p -eq $null) {
        Write-Output ("NOT FOUND: " + $needle)
        return
    }
    $r = $p.Range
    $r.InsertXML($xml)
}

# 1) "If using an Xbox One devkit, set the active solution platform to Gaming.Xbox.XboxOne.x64."
#    -> split around "Gaming.Xbox.XboxOne.x" with gramStart/gramEnd proofErr markers.
Replace-ParaXml "Xbox One devkit" (
  '<w:p w14:paraId="41C15625" w14:textId="77777777" w:rsidR="00B852DF" w:rsidRDefault="00B852DF" w:rsidP="00B852DF">' +
    '<w:r><w:t xml:space="preserve">If using an Xbox One devkit, set the active solution platform to </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Gaming.Xbox.XboxOne.x</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>64.</w:t></w:r>' +
  '</w:p>'
)

# 2) ", set the active solution platform to Gaming.Xbox.Scarlett.x64."
#    -> split around "Gaming.Xbox.Scarlett.x" with gramStart/gramEnd proofErr markers.
Replace-ParaXml "Xbox Series X|S devkit" (
  '<w:p w14:paraId="300C757D" w14:textId="0C9D5FD3" w:rsidR="00B852DF" w:rsidRDefault="00B852DF" w:rsidP="00B852DF">' +
    '<w:r><w:t xml:space="preserve">If using </w:t></w:r>' +
    '<w:r w:rsidR="0083322E"><w:t>an Xbox Series X|S devkit</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, set the active solution platform to </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Gaming.Xbox.Scarlett.x</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>64.</w:t></w:r>' +
  '</w:p>'
)

# 3) "...it takes two calls to the API to get all of the achievements.  "
#    -> split around "all of" with gramStart/gramEnd proofErr markers.
Replace-ParaXml "takes two calls to the API" (
  '<w:p w14:paraId="6A9E6703" w14:textId="4E57867B" w:rsidR="000251A1" w:rsidRDefault="000251A1" w:rsidP="000251A1">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="15"/></w:numPr><w:spacing w:after="160" w:line="259" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Note that because the sample artificially limits the &#8220;take&#8221; count to 1 and there are two achievements registered to this sample, it takes two calls to the API to get </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>all of</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> the achievements.  </w:t></w:r>' +
  '</w:p>'
)

# 4) "Use this to check the status...do not wish to query all the achievements for a title."
#    -> split around "query" with gramStart/gramEnd proofErr markers (the trailing "all" run
#       and the final run are untouched).
Replace-ParaXml "do not wish to query" (
  '<w:p w14:paraId="3C8EE61E" w14:textId="0BE1A209" w:rsidR="000251A1" w:rsidRDefault="000251A1" w:rsidP="000251A1">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="15"/></w:numPr><w:spacing w:after="160" w:line="259" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Use this to check the status of a specific achievement.  You may prefer this if you do not wish to </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>query</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidR="00F635D1"><w:t>all</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> the achievements for a title.</w:t></w:r>' +
  '</w:p>'
)

# 5) "...or use the XblPlayerDataReset tool located here:"
#    -> wrap the bolded "XblPlayerDataReset" run with spellStart/spellEnd proofErr markers
#       (no text change).
Replace-ParaXml "XblPlayerDataReset" (
  '<w:p w14:paraId="10E6D4D6" w14:textId="4C63A92F" w:rsidR="006B1A4D" w:rsidRDefault="00085417" w:rsidP="003D3EF7">' +
    '<w:r><w:t>Once all achievements have been attained by the user, you cannot retrigger the achievements.  You will need to select a new user to show it again</w:t></w:r>' +
    '<w:r w:rsidR="00A45265"><w:t xml:space="preserve">, or use the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r w:rsidR="00A45265" w:rsidRPr="006B1A4D"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>XblPlayerDataReset</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidR="00A45265"><w:t xml:space="preserve"> tool</w:t></w:r>' +
    '<w:r w:rsidR="006B1A4D"><w:t xml:space="preserve"> located here:</w:t></w:r>' +
  '</w:p>'
)
